$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7505.7144
$ws.Range("I32").Value = 6700.5
$ws.Range("J32").Value = 7827.8
$ws.Range("K32").Value = 6700.5
$ws.Range("L32").Value = 7827.8
$ws.Range("M32").Value = -6374.5
$ws.Range("N32").Value = -8479.799999999999

$ws.Range("H38").Value = 4312.1924
$ws.Range("J38").Value = 5683.9287
$ws.Range("L38").Value = 17051.7861
$ws.Range("N38").Value = -17795.7861

$ws.Range("H74").Value = 5832.6
$ws.Range("I74").Value = 5450
$ws.Range("K74").Value = 5450
$ws.Range("M74").Value = -4514

$ws.Range("H77").Value = 5832.6
$ws.Range("I77").Value = 5450
$ws.Range("K77").Value = 27250
$ws.Range("M77").Value = -22570

$ws.Range("H111").Value = 3917.077
$ws.Range("I111").Value = 4211.1113
$ws.Range("J111").Value = 3255.5
$ws.Range("K111").Value = 12633.3339
$ws.Range("L111").Value = 9766.5
$ws.Range("M111").Value = -9566.333899999998
$ws.Range("N111").Value = -15900.5

$ws.Range("H132").Value = 2368.5908
$ws.Range("I132").Value = 2498.0244
$ws.Range("K132").Value = 7494.073199999999
$ws.Range("M132").Value = -4964.073199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3426.8948
$ws.Range("I2").Value = 2756.8125
$ws.Range("K2").Value = 2756.8125
$ws.Range("M2").Value = -2643.8125

$ws.Range("H23").Value = 17000
$ws.Range("I23").Value = 17000
$ws.Range("K23").Value = 17000
$ws.Range("M23").Value = -16741

$ws.Range("H61").Value = 2559.1724
$ws.Range("I61").Value = 2660
$ws.Range("K61").Value = 2660
$ws.Range("M61").Value = -2448

$ws.Range("H74").Value = 58560.906
$ws.Range("I74").Value = 61288.25
$ws.Range("K74").Value = 61288.25
$ws.Range("M74").Value = -60414.25

$ws.Range("H77").Value = 58560.906
$ws.Range("I77").Value = 61288.25
$ws.Range("K77").Value = 306441.25
$ws.Range("M77").Value = -302073.25

$ws.Range("H102").Value = 7681.8
$ws.Range("I102").Value = 7602.25
$ws.Range("K102").Value = 7602.25
$ws.Range("M102").Value = -5980.25

$ws.Range("H116").Value = 3426.8948
$ws.Range("I116").Value = 2756.8125
$ws.Range("K116").Value = 2756.8125
$ws.Range("M116").Value = -462.8125

$ws.Range("H136").Value = 2559.1724
$ws.Range("I136").Value = 2660
$ws.Range("K136").Value = 7980
$ws.Range("M136").Value = -5430

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3426.8948
$ws.Range("I3").Value = 2756.8125
$ws.Range("K3").Value = 2756.8125
$ws.Range("M3").Value = -2642.8125

$ws.Range("H99").Value = 4094.9167
$ws.Range("I99").Value = 2875.4285
$ws.Range("J99").Value = 5802.2
$ws.Range("K99").Value = 2875.4285
$ws.Range("L99").Value = 5802.2
$ws.Range("M99").Value = -1377.4285
$ws.Range("N99").Value = -8798.200000000001

$ws.Range("H105").Value = 10183.459
$ws.Range("I105").Value = 24804.555
$ws.Range("K105").Value = 24804.555
$ws.Range("M105").Value = -23057.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3093.1155
$ws.Range("I31").Value = 2032.5454
$ws.Range("J31").Value = 3870.8667
$ws.Range("K31").Value = 2032.5454
$ws.Range("L31").Value = 3870.8667
$ws.Range("M31").Value = -1737.5454
$ws.Range("N31").Value = -4460.8667

$ws.Range("H34").Value = 3093.1155
$ws.Range("I34").Value = 2032.5454
$ws.Range("J34").Value = 3870.8667
$ws.Range("K34").Value = 2032.5454
$ws.Range("L34").Value = 3870.8667
$ws.Range("M34").Value = -1830.5454
$ws.Range("N34").Value = -4274.8667

$ws.Range("H107").Value = 524.5
$ws.Range("I107").Value = 466.33334
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 466.33334
$ws.Range("L107").Value = 699
$ws.Range("M107").Value = 1453.66666
$ws.Range("N107").Value = -4539

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 7059
$ws.Range("J44").Value = 3041.6667
$ws.Range("L44").Value = 9125.000100000001
$ws.Range("N44").Value = -9921.000100000001

$ws.Range("H113").Value = 1096.6875
$ws.Range("I113").Value = 416.3846
$ws.Range("J113").Value = 4044.6667
$ws.Range("K113").Value = 1249.1538
$ws.Range("L113").Value = 12134.0001
$ws.Range("M113").Value = 920.8462
$ws.Range("N113").Value = -16474.0001

$ws.Range("H122").Value = 3048.7
$ws.Range("I122").Value = 847
$ws.Range("J122").Value = 3599.125
$ws.Range("K122").Value = 7623
$ws.Range("L122").Value = 32392.125
$ws.Range("M122").Value = -5173
$ws.Range("N122").Value = -37292.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4097.6
$ws.Range("J113").Value = 4249.5
$ws.Range("L113").Value = 4249.5
$ws.Range("N113").Value = -8589.5

$ws.Range("H123").Value = 50566.285
$ws.Range("J123").Value = 47991.25
$ws.Range("L123").Value = 47991.25
$ws.Range("N123").Value = -52891.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4021
$ws.Range("I7").Value = 3401.077
$ws.Range("J7").Value = 6035.75
$ws.Range("K7").Value = 3401.077
$ws.Range("L7").Value = 6035.75
$ws.Range("M7").Value = -3289.077
$ws.Range("N7").Value = -6259.75

$ws.Range("H16").Value = 1589.3823
$ws.Range("I16").Value = 1205.7391
$ws.Range("K16").Value = 1205.7391
$ws.Range("M16").Value = -1035.7391

$ws.Range("H61").Value = 2439.9048
$ws.Range("I61").Value = 2234.1177
$ws.Range("J61").Value = 3314.5
$ws.Range("K61").Value = 2234.1177
$ws.Range("L61").Value = 3314.5
$ws.Range("M61").Value = -2032.1177
$ws.Range("N61").Value = -3718.5

$ws.Range("H69").Value = 50000
$ws.Range("I69").Value = 50000
$ws.Range("K69").Value = 50000
$ws.Range("M69").Value = -49189

$ws.Range("H72").Value = 50000
$ws.Range("I72").Value = 50000
$ws.Range("K72").Value = 150000
$ws.Range("M72").Value = -145944

$ws.Range("H113").Value = 2439.9048
$ws.Range("I113").Value = 2234.1177
$ws.Range("J113").Value = 3314.5
$ws.Range("K113").Value = 2234.1177
$ws.Range("L113").Value = 3314.5
$ws.Range("M113").Value = -64.11769999999979
$ws.Range("N113").Value = -7654.5

$ws.Range("H126").Value = 4021
$ws.Range("I126").Value = 3401.077
$ws.Range("J126").Value = 6035.75
$ws.Range("K126").Value = 10203.231
$ws.Range("L126").Value = 18107.25
$ws.Range("M126").Value = -7733.231
$ws.Range("N126").Value = -23047.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 74994.5
$ws.Range("J119").Value = 74994.5
$ws.Range("L119").Value = 74994.5
$ws.Range("N119").Value = -84670.5
